# Auto-generated edit script: update column F ('想去人数' / interest counts)
# across sheets 1 (展览), 2 (演出), and 4 (全部类型) to match refreshed scrape data.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F5").Value = 1146
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 145
$ws.Range("F8").Value = 239
$ws.Range("F9").Value = 7115
$ws.Range("F12").Value = 5437
$ws.Range("F13").Value = 10
$ws.Range("F14").Value = 71
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F20").Value = 146
$ws.Range("F22").Value = 161
$ws.Range("F23").Value = 104
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 1870
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 0
$ws.Range("F34").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 980
$ws.Range("F47").Value = 0
$ws.Range("F49").Value = 0

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 192
$ws.Range("F13").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("F15").Value = 97
$ws.Range("F17").Value = 0
$ws.Range("F19").Value = 7
$ws.Range("F21").Value = 2

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 1146
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 7115
$ws.Range("F19").Value = 6189
$ws.Range("F20").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 146
$ws.Range("F25").Value = 0
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 10047
$ws.Range("F32").Value = 2082
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 5178
$ws.Range("F41").Value = 1196
$ws.Range("F42").Value = 652
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 168
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 1359

$wb.Save()
